$wb = $excel.ActiveWorkbook

# --- caseType1: move the "Label" column (E) to just after "Case Property" (A),
# i.e. make it the new column B. Group/Data Type/Description each slide one
# column to the right (B->C, C->D, D->E); "Deprecated" in F is untouched. ---
$ws1 = $wb.Worksheets.Item("caseType1")
$ws1.Columns.Item(5).Cut()
$ws1.Columns.Item(2).Insert()
$ws1.Columns.Item(3).ColumnWidth = 24.1666667
$ws1.Columns.Item(4).ColumnWidth = 25.3333333
$ws1.Columns.Item(5).ColumnWidth = 9.8333333

# --- caseType2: same column reorder ---
$ws2 = $wb.Worksheets.Item("caseType2")
$ws2.Columns.Item(5).Cut()
$ws2.Columns.Item(2).Insert()

# Restore the on-screen selection: both sheets end up with F1 selected,
# but caseType1 is the active/visible tab (selected last).
$ws2.Range("F1").Select()
$ws1.Range("F1").Select()
